$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 10.75
$ws.Range("I5").Value = 10.75
$ws.Range("K5").Value = 10.75
$ws.Range("M5").Value = 104.25

$ws.Range("H18").Value = 1094.7931
$ws.Range("I18").Value = 1080.3214
$ws.Range("J18").Value = 1500
$ws.Range("K18").Value = 1080.3214
$ws.Range("L18").Value = 1500
$ws.Range("M18").Value = -796.3214
$ws.Range("N18").Value = -2068

$ws.Range("H132").Value = 3249080.2
$ws.Range("I132").Value = 4084138.8
$ws.Range("J132").Value = 1630.7778
$ws.Range("K132").Value = 12252416.4
$ws.Range("L132").Value = 4892.3334
$ws.Range("M132").Value = -12249886.4
$ws.Range("N132").Value = -9952.3334

$ws.Range("H137").Value = 1203.4108
$ws.Range("I137").Value = 624
$ws.Range("J137").Value = 2941.6428
$ws.Range("K137").Value = 1872
$ws.Range("L137").Value = 8824.928400000001
$ws.Range("M137").Value = 678
$ws.Range("N137").Value = -13924.9284

$ws.Range("H141").Value = 2681.8333
$ws.Range("I141").Value = 1644.75
$ws.Range("J141").Value = 3624.6365
$ws.Range("K141").Value = 4934.25
$ws.Range("L141").Value = 10873.9095
$ws.Range("M141").Value = 245.75
$ws.Range("N141").Value = -21233.9095

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20408882
$ws.Range("I2").Value = 28571916
$ws.Range("J2").Value = 1291.5
$ws.Range("K2").Value = 28571916
$ws.Range("L2").Value = 1291.5
$ws.Range("M2").Value = -28571803
$ws.Range("N2").Value = -1517.5

$ws.Range("H32").Value = 1709.86
$ws.Range("I32").Value = 1709.86
$ws.Range("K32").Value = 1709.86
$ws.Range("M32").Value = -1422.86

$ws.Range("H61").Value = 2792.6667
$ws.Range("I61").Value = 1918.6666
$ws.Range("J61").Value = 3666.6667
$ws.Range("K61").Value = 1918.6666
$ws.Range("L61").Value = 3666.6667
$ws.Range("M61").Value = -1706.6666
$ws.Range("N61").Value = -4090.6667

$ws.Range("H74").Value = 1149.3334
$ws.Range("I74").Value = 1078.4706
$ws.Range("J74").Value = 1321.4286
$ws.Range("K74").Value = 1078.4706
$ws.Range("L74").Value = 1321.4286
$ws.Range("M74").Value = -204.4706000000001
$ws.Range("N74").Value = -3069.4286

$ws.Range("H77").Value = 1149.3334
$ws.Range("I77").Value = 1078.4706
$ws.Range("J77").Value = 1321.4286
$ws.Range("K77").Value = 5392.353000000001
$ws.Range("L77").Value = 6607.143
$ws.Range("M77").Value = -1024.353000000001
$ws.Range("N77").Value = -15343.143

$ws.Range("H116").Value = 20408882
$ws.Range("I116").Value = 28571916
$ws.Range("J116").Value = 1291.5
$ws.Range("K116").Value = 28571916
$ws.Range("L116").Value = 1291.5
$ws.Range("M116").Value = -28569622
$ws.Range("N116").Value = -5879.5

$ws.Range("H122").Value = 1832.091
$ws.Range("I122").Value = 1828.9412
$ws.Range("K122").Value = 5486.8236
$ws.Range("M122").Value = -3036.8236

$ws.Range("H136").Value = 2792.6667
$ws.Range("I136").Value = 1918.6666
$ws.Range("J136").Value = 3666.6667
$ws.Range("K136").Value = 5755.9998
$ws.Range("L136").Value = 11000.0001
$ws.Range("M136").Value = -3205.9998
$ws.Range("N136").Value = -16100.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20408882
$ws.Range("I3").Value = 28571916
$ws.Range("J3").Value = 1291.5
$ws.Range("K3").Value = 28571916
$ws.Range("L3").Value = 1291.5
$ws.Range("M3").Value = -28571802
$ws.Range("N3").Value = -1519.5

$ws.Range("H107").Value = 1133.878
$ws.Range("I107").Value = 1123.3667
$ws.Range("J107").Value = 1162.5454
$ws.Range("K107").Value = 1123.3667
$ws.Range("L107").Value = 1162.5454
$ws.Range("M107").Value = 796.6333
$ws.Range("N107").Value = -5002.5454

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3549287.5
$ws.Range("I31").Value = 2485.2778
$ws.Range("J31").Value = 15157004
$ws.Range("K31").Value = 2485.2778
$ws.Range("L31").Value = 15157004
$ws.Range("M31").Value = -2190.2778
$ws.Range("N31").Value = -15157594

$ws.Range("H34").Value = 3549287.5
$ws.Range("I34").Value = 2485.2778
$ws.Range("J34").Value = 15157004
$ws.Range("K34").Value = 2485.2778
$ws.Range("L34").Value = 15157004
$ws.Range("M34").Value = -2283.2778
$ws.Range("N34").Value = -15157408

$ws.Range("H58").Value = 899.1111
$ws.Range("I58").Value = 911.5
$ws.Range("K58").Value = 911.5
$ws.Range("M58").Value = -708.5

$ws.Range("H132").Value = 5210372
$ws.Range("I132").Value = 938.5454999999999
$ws.Range("J132").Value = 9618354
$ws.Range("K132").Value = 2815.6365
$ws.Range("L132").Value = 28855062
$ws.Range("M132").Value = -285.6364999999996
$ws.Range("N132").Value = -28860122

$ws.Range("H134").Value = 1055
$ws.Range("I134").Value = 1055
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3165
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -630
$ws.Range("N134").ClearContents()

$ws.Range("H136").Value = 899.1111
$ws.Range("I136").Value = 911.5
$ws.Range("K136").Value = 2734.5
$ws.Range("M136").Value = -184.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 394
$ws.Range("I4").Value = 131.8
$ws.Range("J4").Value = 1049.5
$ws.Range("K4").Value = 395.4
$ws.Range("L4").Value = 3148.5
$ws.Range("M4").Value = -283.4
$ws.Range("N4").Value = -3372.5

$ws.Range("H26").Value = 12500729
$ws.Range("I26").Value = 45.333332
$ws.Range("J26").Value = 20001138
$ws.Range("K26").Value = 135.999996
$ws.Range("L26").Value = 60003414
$ws.Range("M26").Value = 152.000004
$ws.Range("N26").Value = -60003990

$ws.Range("H107").Value = 437.7619
$ws.Range("I107").Value = 847.5
$ws.Range("J107").Value = 341.35294
$ws.Range("K107").Value = 2542.5
$ws.Range("L107").Value = 1024.05882
$ws.Range("M107").Value = -622.5
$ws.Range("N107").Value = -4864.05882

$ws.Range("H113").Value = 568.43054
$ws.Range("J113").Value = 622.875
$ws.Range("L113").Value = 1868.625
$ws.Range("N113").Value = -6208.625

$ws.Range("H131").Value = 1209716.1
$ws.Range("I131").Value = 11343
$ws.Range("J131").Value = 1355859.1
$ws.Range("K131").Value = 34029
$ws.Range("L131").Value = 4067577.3
$ws.Range("M131").Value = -28989
$ws.Range("N131").Value = -4077657.3

$ws.Range("H141").Value = 3057.08
$ws.Range("I141").Value = 2717.2104
$ws.Range("J141").Value = 4133.3335
$ws.Range("K141").Value = 8151.6312
$ws.Range("L141").Value = 12400.0005
$ws.Range("M141").Value = -2971.6312
$ws.Range("N141").Value = -22760.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1738.1666
$ws.Range("I126").Value = 1800.421
$ws.Range("J126").Value = 1501.6
$ws.Range("K126").Value = 5401.263
$ws.Range("L126").Value = 4504.799999999999
$ws.Range("M126").Value = -2931.263
$ws.Range("N126").Value = -9444.799999999999

$ws.Range("H132").Value = 60776
$ws.Range("I132").Value = 85029.164
$ws.Range("J132").Value = 2568.4
$ws.Range("K132").Value = 255087.492
$ws.Range("L132").Value = 7705.200000000001
$ws.Range("M132").Value = -252557.492
$ws.Range("N132").Value = -12765.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2202.0952
$ws.Range("I40").Value = 2212.2
$ws.Range("K40").Value = 2212.2
$ws.Range("M40").Value = -2076.2

$ws.Range("H122").Value = 3587.8333
$ws.Range("I122").Value = 4509
$ws.Range("J122").Value = 2666.6667
$ws.Range("K122").Value = 13527
$ws.Range("L122").Value = 8000.000100000001
$ws.Range("M122").Value = -11077
$ws.Range("N122").Value = -12900.0001

$ws.Range("H132").Value = 6405.4634
$ws.Range("I132").Value = 8734.666999999999
$ws.Range("K132").Value = 26204.001
$ws.Range("M132").Value = -23674.001

$ws.Range("H133").Value = 26000
$ws.Range("J133").Value = 26000
$ws.Range("L133").Value = 26000
$ws.Range("N133").Value = -31060

$ws.Range("H136").Value = 5045.057
$ws.Range("I136").Value = 5823.08
$ws.Range("J136").Value = 3100
$ws.Range("K136").Value = 17469.24
$ws.Range("L136").Value = 9300
$ws.Range("M136").Value = -14919.24
$ws.Range("N136").Value = -14400

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3611.353
$ws.Range("I122").Value = 4350.5
$ws.Range("K122").Value = 13051.5
$ws.Range("M122").Value = -10601.5

$ws.Range("H132").Value = 1343.3922
$ws.Range("I132").Value = 1255.3673
$ws.Range("J132").Value = 3500
$ws.Range("K132").Value = 3766.1019
$ws.Range("L132").Value = 10500
$ws.Range("M132").Value = -1236.1019
$ws.Range("N132").Value = -15560

$ws.Range("H136").Value = 1148.3455
$ws.Range("I136").Value = 1092.3077
$ws.Range("K136").Value = 3276.9231
$ws.Range("M136").Value = -726.9231
